$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.193.24'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '3.087.18'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.28'
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.99'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '3.081.95'
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.507'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.14'
$ws.Range("E11").Value = '  -6.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.472'
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.11'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '3.602.31'
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").Value = '64.262.61'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").Value = '3.097.81'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.74'
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '481.48'
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.03'
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.677'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.54'
$ws.Range("E23").Value = '  +3.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.11'
$ws.Range("E24").Value = '  +10.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.40'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.00'
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.05'
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.32'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.65'
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.20'
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.39'
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.97'
$ws.Range("E37").Value = '  +14.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0408'
$ws.Range("E38").Value = '  +2.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '442.03'
$ws.Range("E39").Value = '  -5.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0812'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").Value = '2.955.83'
$ws.Range("E41").Value = '  -3.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.21'
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("E43").Value = '  -4.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.15'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  +3.50%  '
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.37'
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("D50").Value = '0.0₃0516'
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.09'
$ws.Range("E51").Value = '  -0.12%  '
